# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update the "last updated" timestamp (A1 / shared string) ---
$ws.Range("A1").Value = "Datos actualizados a 24 de Mayo de 2020 a las 19:05"

# --- Swap the Mexico / Chile rows (row 19 was Mexico, row 20 was Chile) ---
# After the update Chile moves to row 19 (with new figures) and Mexico moves
# to row 20 keeping the figures that used to belong to the row-19 country.
$ws.Range("A19").Value = "Chile"
$ws.Range("A20").Value = "Mexico"

# --- Numeric data updates ---

# Row 4: Estados Unidos
$ws.Cells.Item(4, 2).Value = 1673806
$ws.Cells.Item(4, 3).Value = 6978
$ws.Cells.Item(4, 4).Value = 448994
$ws.Cells.Item(4, 5).Value = 1125981
$ws.Cells.Item(4, 7).Value = 148
$ws.Cells.Item(4, 8).Value = 98831

# Row 11: Alemania
$ws.Cells.Item(11, 2).Value = 180105
$ws.Cells.Item(11, 3).Value = 119
$ws.Cells.Item(11, 5).Value = 11434

# Row 13: India
$ws.Cells.Item(13, 2).Value = 138041
$ws.Cells.Item(13, 3).Value = 6618
$ws.Cells.Item(13, 4).Value = 57429
$ws.Cells.Item(13, 5).Value = 76598
$ws.Cells.Item(13, 7).Value = 146
$ws.Cells.Item(13, 8).Value = 4014

# Row 16: Canada
$ws.Cells.Item(16, 2).Value = 84082
$ws.Cells.Item(16, 3).Value = 461
$ws.Cells.Item(16, 4).Value = 43640
$ws.Cells.Item(16, 5).Value = 34062

# Row 19: now Chile (new, higher figures)
$ws.Cells.Item(19, 2).Value = 69102
$ws.Cells.Item(19, 3).Value = 3709
$ws.Cells.Item(19, 4).Value = 28148
$ws.Cells.Item(19, 5).Value = 40236
$ws.Cells.Item(19, 7).Value = 45
$ws.Cells.Item(19, 8).Value = 718

# Row 20: now Mexico (figures that used to be on row 19)
$ws.Cells.Item(20, 2).Value = 65856
$ws.Cells.Item(20, 3).Value = 3329
$ws.Cells.Item(20, 4).Value = 44919
$ws.Cells.Item(20, 5).Value = 13758
$ws.Cells.Item(20, 7).Value = 190
$ws.Cells.Item(20, 8).Value = 7179

# Row 45: Republica Dominicana
$ws.Cells.Item(45, 2).Value = 14801
$ws.Cells.Item(45, 3).Value = 379
$ws.Cells.Item(45, 4).Value = 8133
$ws.Cells.Item(45, 5).Value = 6210

# Row 60: Marruecos
$ws.Cells.Item(60, 2).Value = 7433
$ws.Cells.Item(60, 3).Value = 27
$ws.Cells.Item(60, 4).Value = 4703
$ws.Cells.Item(60, 5).Value = 2531
$ws.Cells.Item(60, 7).Value = 1
$ws.Cells.Item(60, 8).Value = 199

# Row 118: Paraguay
$ws.Cells.Item(118, 2).Value = 862
$ws.Cells.Item(118, 3).Value = 12
$ws.Cells.Item(118, 4).Value = 307
$ws.Cells.Item(118, 5).Value = 544

$wb.Save()
